$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking data columns (totalRuns/totalBalls/total4s/
# total6s/sr) to be stored as text so values like "0", "1", "33.33" keep
# their original text representation instead of being silently
# reinterpreted as numbers by Excel's automatic type detection.
$ws.Range("G2:K4").NumberFormat = "@"

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "venue"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "result"
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"
$ws.Range("F1").Value = "batsman"
$ws.Range("G1").Value = "totalRuns"
$ws.Range("H1").Value = "totalBalls"
$ws.Range("I1").Value = "total4s"
$ws.Range("J1").Value = "total6s"
$ws.Range("K1").Value = "sr"

# ---- Row 2 (existing match, now with ownTeam/oppTeam inserted) ----
$ws.Range("A2").Value = " Abu Dhabi"
$ws.Range("B2").Value = " October 07 2020"
$ws.Range("C2").Value = "KKR won by 10 runs"
$ws.Range("D2").Value = "Kolkata Knight Riders"
$ws.Range("E2").Value = "Chennai Super Kings"
$ws.Range("F2").Value = "Shivam Mavi "
$ws.Range("G2").Value = "0"
$ws.Range("H2").Value = "1"
$ws.Range("I2").Value = "0"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "0.00"

# ---- Row 3 (new match) ----
$ws.Range("A3").Value = " Sharjah"
$ws.Range("B3").Value = " October 03 2020"
$ws.Range("C3").Value = "Capitals won by 18 runs"
$ws.Range("D3").Value = "Kolkata Knight Riders"
$ws.Range("E3").Value = "Delhi Capitals"
$ws.Range("F3").Value = "Shivam Mavi "
$ws.Range("G3").Value = "1"
$ws.Range("H3").Value = "3"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "33.33"

# ---- Row 4 (new match) ----
$ws.Range("A4").Value = " Abu Dhabi"
$ws.Range("B4").Value = " September 23 2020"
$ws.Range("C4").Value = "Mumbai won by 49 runs"
$ws.Range("D4").Value = "Kolkata Knight Riders"
$ws.Range("E4").Value = "Mumbai Indians"
$ws.Range("F4").Value = "Shivam Mavi "
$ws.Range("G4").Value = "9"
$ws.Range("H4").Value = "10"
$ws.Range("I4").Value = "1"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "90.00"
